# Applies the "Correção de pequeno erro" commit:
#  1. Insert a new "IsolatedDeps" column between "TestIWithDeps" (F) and "Precision" (old G),
#     pushing Precision..F2Deps from G:L to H:M.
#  2. Fill the new IsolatedDeps column (G2:G14) with a single space, except row 6 which gets "[]".
#  3. Fix the changed-files list order in F6 (TestIWithDeps for task 126).
#  4. Replace the literal "Null" placeholders (now in J and M after the shift) with "[]".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert new column at G - shifts Precision..F2Deps (G:L) to H:M
$ws.Columns("G").Insert()

# 2. New column header + values
$ws.Range("G1").Value = "IsolatedDeps"

$ws.Range("G2").Value  = " "
$ws.Range("G3").Value  = " "
$ws.Range("G4").Value  = " "
$ws.Range("G5").Value  = " "
$ws.Range("G6").Value  = "[]"
$ws.Range("G7").Value  = " "
$ws.Range("G8").Value  = " "
$ws.Range("G9").Value  = " "
$ws.Range("G10").Value = " "
$ws.Range("G11").Value = " "
$ws.Range("G12").Value = " "
$ws.Range("G13").Value = " "
$ws.Range("G14").Value = " "

# 3. Correct the ordering of the changed-files list for row 6 (task 126)
$ws.Range("F6").Value = "['app/views/abingo_dashboard/_experiment_row.html.haml', 'app/models/question.rb', 'app/views/abingo_dashboard/index.html.haml']"

# 4. Replace "Null" placeholders with "[]" (columns J and M, formerly I and L)
$ws.Range("J3").Value = "[]"
$ws.Range("M3").Value = "[]"
$ws.Range("J4").Value = "[]"
$ws.Range("M4").Value = "[]"
$ws.Range("J5").Value = "[]"
$ws.Range("M5").Value = "[]"
$ws.Range("J6").Value = "[]"
$ws.Range("M6").Value = "[]"
$ws.Range("J8").Value = "[]"
$ws.Range("M8").Value = "[]"
$ws.Range("J9").Value = "[]"
$ws.Range("M9").Value = "[]"
$ws.Range("J10").Value = "[]"
$ws.Range("M10").Value = "[]"
$ws.Range("J11").Value = "[]"
$ws.Range("M11").Value = "[]"
